# Add the new "2018-07-25" experiment tab (3rd protein-unfolding run), fill in
# its data/formulas, add the corresponding scatter chart with a linear
# trendline, and mirror the accompanying selection / active-sheet changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New worksheet, appended after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2018-07-25"

# ---------------------------------------------------------------------------
# 2. Headers (reuse the shared strings already used on the other tabs).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Sample"
$ws.Range("B1").Value = "R1"
$ws.Range("C1").Value = "R2"
$ws.Range("D1").Value = "Avg"
$ws.Range("E1").Value = "Normalized"
$ws.Range("F1").Value = "Predicted"

# ---------------------------------------------------------------------------
# 3. Data rows (R1/R2 replicate readings, averaged, normalized against the
#    blank (row 6), then converted with the run's calibration line).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 2.356
$ws.Range("C2").Value = 1.736

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1.166
$ws.Range("C3").Value = 1.211

$ws.Range("A4").Value = 0.5
$ws.Range("B4").Value = 0.681
$ws.Range("C4").Value = 0.669

$ws.Range("A5").Value = 0.25
$ws.Range("B5").Value = 0.39
$ws.Range("C5").Value = 0.374

$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 0.068
$ws.Range("C6").Value = 0.067
# Row 6 "Avg" is the literal blank baseline (not a formula average).
$ws.Range("D6").Value = 0.072

$ws.Range("A7").Value = "samp"
$ws.Range("B7").Value = 0.767
$ws.Range("C7").Value = 0.782

# Avg = AVERAGE(B:C) for every row except the literal blank row (6).
$ws.Range("D2").Formula = "=AVERAGE(B2:C2)"
$ws.Range("D3").Formula = "=AVERAGE(B3:C3)"
$ws.Range("D4").Formula = "=AVERAGE(B4:C4)"
$ws.Range("D5").Formula = "=AVERAGE(B5:C5)"
$ws.Range("D7").Formula = "=AVERAGE(B7:C7)"

# Normalized = Avg - blank Avg, for every data row (including the blank itself).
$ws.Range("E2").Formula = "=D2-`$D`$6"
$ws.Range("E3").Formula = "=D3-`$D`$6"
$ws.Range("E4").Formula = "=D4-`$D`$6"
$ws.Range("E5").Formula = "=D5-`$D`$6"
$ws.Range("E6").Formula = "=D6-`$D`$6"
$ws.Range("E7").Formula = "=D7-`$D`$6"

# Predicted = calibration line applied to the normalized value.
$ws.Range("F2").Formula = "=(1.0182*E2) - 0.0652"
$ws.Range("F3").Formula = "=(1.0182*E3) - 0.0652"
$ws.Range("F4").Formula = "=(1.0182*E4) - 0.0652"
$ws.Range("F5").Formula = "=(1.0182*E5) - 0.0652"
$ws.Range("F6").Formula = "=(1.0182*E6) - 0.0652"
$ws.Range("F7").Formula = "=(1.0182*E7) - 0.0652"

# ---------------------------------------------------------------------------
# 4. Column widths (match the bestFit widths used on the new tab).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7.875
$ws.Columns.Item(2).ColumnWidth = 5.875
$ws.Columns.Item(3).ColumnWidth = 5.875
$ws.Columns.Item(4).ColumnWidth = 6.875
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 9.875

# ---------------------------------------------------------------------------
# 5. Scatter chart of Normalized (x) vs. Sample amount (y) with a linear
#    trendline + equation, same style as the other tabs' charts.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Add(410, 130, 410, 227)
$co.Chart.ChartType = 74
$co.Chart.SeriesCollection().NewSeries()
$ser = $co.Chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,'2018-07-25'!`$E`$2:`$E`$6,'2018-07-25'!`$A`$2:`$A`$6,1)"

$trend = $ser.Trendlines().Add()
$trend.Type = -4132
$trend.DisplayEquation = $true
$trend.DisplayRSquared = $false

# ---------------------------------------------------------------------------
# 6. Selection / active-sheet bookkeeping to match the edited workbook.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:G9").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C12").Select()

$ws.Range("G7").Select()
$ws.Activate()
